$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")

# Update the note in L44: append the "UPDATE" clarification.
# This removes the now-unused old shared string "Does not match what is
# published" and appends the new text to the shared-string table.
$ws.Range("L44").Value = "Does not match what is published: UPDATE: Because github model is close but not quite as good as what went in the paper."

# New row 46: DenseDepth (Wasserstein histogram matching) label only.
# (Set before row 41 below so the shared-string table gets the
# DenseDepth/DORN Wasserstein labels in that order, matching the source file.)
$ws.Range("A46").Value = "DenseDepth (Wasserstein histogram matching)"

# New row 41: DORN (Wasserstein histogram matching) results, inserted into
# the previously-empty row right after row 40 (DORN GT histogram matching).
$ws.Range("A41").Value = "DORN (Wasserstein histogram matching)"
$ws.Range("B41").Value = 0.847427449419342
$ws.Range("C41").Value = 0.95332383895321304
$ws.Range("D41").Value = 0.982672920285379
$ws.Range("F41").Value = 0.49932645306856899
$ws.Range("G41").Value = 0.117189220622728
$ws.Range("I41").Value = 0.053489108434636203

# Make row 44 a bit taller to fit the longer note.
$ws.Rows.Item(44).RowHeight = 31

# Update the view: active selection (moved down as more rows were filled in).
$ws.Activate()
$ws.Range("J45").Select()
